$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H93").Value = 39286.957
$ws.Range("J93").Value = 39286.957
$ws.Range("L93").Value = 39286.957
$ws.Range("N93").Value = -44278.957
$ws.Range("H112").Value = 1354.1464
$ws.Range("J112").Value = 1368
$ws.Range("L112").Value = 4104
$ws.Range("N112").Value = -6320
$ws.Range("H129").Value = 1530.9636
$ws.Range("J129").Value = 1585.6923
$ws.Range("L129").Value = 4757.0769
$ws.Range("N129").Value = -14757.0769
$ws.Range("H133").Value = 54714.285
$ws.Range("J133").Value = 54714.285
$ws.Range("L133").Value = 54714.285
$ws.Range("N133").Value = -64834.285
$ws.Range("H136").Value = 52872.5
$ws.Range("J136").Value = 52872.5
$ws.Range("L136").Value = 52872.5
$ws.Range("N136").Value = -63072.5
$ws.Range("H137").Value = 664277
$ws.Range("I137").Value = 1645180
$ws.Range("J137").Value = 2737.7441
$ws.Range("K137").Value = 4935540
$ws.Range("L137").Value = 8213.2323
$ws.Range("M137").Value = -4932990
$ws.Range("N137").Value = -13313.2323
$ws.Range("H138").Value = 2996.5386
$ws.Range("I138").Value = 1835
$ws.Range("J138").Value = 3375.8164
$ws.Range("K138").Value = 5505
$ws.Range("L138").Value = 10127.4492
$ws.Range("M138").Value = -365
$ws.Range("N138").Value = -20407.4492

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H32").Value = 4522.0645
$ws.Range("I32").Value = 4270.56
$ws.Range("K32").Value = 4270.56
$ws.Range("M32").Value = -3983.56
$ws.Range("H61").Value = 1395.25
$ws.Range("I61").Value = 1308.8572
$ws.Range("J61").Value = 2000
$ws.Range("K61").Value = 1308.8572
$ws.Range("L61").Value = 2000
$ws.Range("M61").Value = -1096.8572
$ws.Range("N61").Value = -2424
$ws.Range("H103").Value = 34750
$ws.Range("J103").Value = 34750
$ws.Range("L103").Value = 34750
$ws.Range("N103").Value = -37094
$ws.Range("H110").Value = 729.2105
$ws.Range("I110").Value = 565.1429000000001
$ws.Range("J110").Value = 1188.6
$ws.Range("K110").Value = 565.1429000000001
$ws.Range("L110").Value = 1188.6
$ws.Range("M110").Value = 1479.8571
$ws.Range("N110").Value = -5278.6
$ws.Range("H136").Value = 1395.25
$ws.Range("I136").Value = 1308.8572
$ws.Range("J136").Value = 2000
$ws.Range("K136").Value = 3926.5716
$ws.Range("L136").Value = 6000
$ws.Range("M136").Value = -1376.5716
$ws.Range("N136").Value = -11100
$ws.Range("H137").Value = 40466
$ws.Range("J137").Value = 40466
$ws.Range("L137").Value = 40466
$ws.Range("N137").Value = -50666

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H86").Value = 2704.7273
$ws.Range("I86").Value = 2705.2
$ws.Range("J86").Value = 2700
$ws.Range("K86").Value = 2705.2
$ws.Range("L86").Value = 2700
$ws.Range("M86").Value = -1582.2
$ws.Range("N86").Value = -4946
$ws.Range("H89").Value = 2704.7273
$ws.Range("I89").Value = 2705.2
$ws.Range("J89").Value = 2700
$ws.Range("K89").Value = 13526
$ws.Range("L89").Value = 13500
$ws.Range("M89").Value = -7910
$ws.Range("N89").Value = -24732
$ws.Range("H95").Value = 29166.5
$ws.Range("J95").Value = 29166.5
$ws.Range("L95").Value = 29166.5
$ws.Range("N95").Value = -34658.5
$ws.Range("H105").Value = 1881.697
$ws.Range("I105").Value = 1840.9683
$ws.Range("J105").Value = 2737
$ws.Range("K105").Value = 1840.9683
$ws.Range("L105").Value = 2737
$ws.Range("M105").Value = -93.9683
$ws.Range("N105").Value = -6231
$ws.Range("H107").Value = 965.05554
$ws.Range("I107").Value = 850.25
$ws.Range("K107").Value = 850.25
$ws.Range("M107").Value = 1069.75
$ws.Range("H137").Value = 55420.832
$ws.Range("J137").Value = 55420.832
$ws.Range("L137").Value = 55420.832
$ws.Range("N137").Value = -65620.83199999999

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H9").Value = 49800
$ws.Range("J9").Value = 49800
$ws.Range("L9").Value = 49800
$ws.Range("N9").Value = -50136
$ws.Range("H22").Value = 606.8570999999999
$ws.Range("I22").Value = 371
$ws.Range("J22").Value = 1104.7778
$ws.Range("K22").Value = 371
$ws.Range("L22").Value = 1104.7778
$ws.Range("M22").Value = -21
$ws.Range("N22").Value = -1804.7778
$ws.Range("H31").Value = 200596.31
$ws.Range("I31").Value = 466840.1
$ws.Range("J31").Value = 2620.1538
$ws.Range("K31").Value = 466840.1
$ws.Range("L31").Value = 2620.1538
$ws.Range("M31").Value = -466545.1
$ws.Range("N31").Value = -3210.1538
$ws.Range("H34").Value = 200596.31
$ws.Range("I34").Value = 466840.1
$ws.Range("J34").Value = 2620.1538
$ws.Range("K34").Value = 466840.1
$ws.Range("L34").Value = 2620.1538
$ws.Range("M34").Value = -466638.1
$ws.Range("N34").Value = -3024.1538
$ws.Range("H124").Value = 45674.332
$ws.Range("J124").Value = 45674.332
$ws.Range("L124").Value = 45674.332
$ws.Range("N124").Value = -50584.332
$ws.Range("H132").Value = 3582.5
$ws.Range("I132").Value = 2554.7036
$ws.Range("J132").Value = 6665.8887
$ws.Range("K132").Value = 7664.110799999999
$ws.Range("L132").Value = 19997.6661
$ws.Range("M132").Value = -5134.110799999999
$ws.Range("N132").Value = -25057.6661
$ws.Range("H134").Value = 6198.0454
$ws.Range("I134").Value = 6239.8423
$ws.Range("J134").Value = 5933.3335
$ws.Range("K134").Value = 18719.5269
$ws.Range("L134").Value = 17800.0005
$ws.Range("M134").Value = -16184.5269
$ws.Range("N134").Value = -22870.0005

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H68").Value = 2417.075
$ws.Range("I68").Value = 778.70966
$ws.Range("J68").Value = 3453.5918
$ws.Range("K68").Value = 2336.12898
$ws.Range("L68").Value = 10360.7754
$ws.Range("M68").Value = -1525.12898
$ws.Range("N68").Value = -11982.7754
$ws.Range("H71").Value = 2417.075
$ws.Range("I71").Value = 778.70966
$ws.Range("J71").Value = 3453.5918
$ws.Range("K71").Value = 7008.38694
$ws.Range("L71").Value = 31082.3262
$ws.Range("M71").Value = -2952.38694
$ws.Range("N71").Value = -39194.3262
$ws.Range("H113").Value = 3788656.5
$ws.Range("I113").Value = 652.8182
$ws.Range("J113").Value = 11364664
$ws.Range("K113").Value = 1958.4546
$ws.Range("L113").Value = 34093992
$ws.Range("M113").Value = 211.5454
$ws.Range("N113").Value = -34098332
$ws.Range("H132").Value = 2157.238
$ws.Range("J132").Value = 2537.375
$ws.Range("L132").Value = 22836.375
$ws.Range("N132").Value = -27896.375
$ws.Range("H137").Value = 3223.6086
$ws.Range("I137").Value = 3529.077
$ws.Range("J137").Value = 2826.5
$ws.Range("K137").Value = 10587.231
$ws.Range("L137").Value = 8479.5
$ws.Range("M137").Value = -5487.231
$ws.Range("N137").Value = -18679.5

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H46").Value = 23337.625
$ws.Range("J46").Value = 24225.143
$ws.Range("L46").Value = 24225.143
$ws.Range("N46").Value = -24537.143
$ws.Range("H70").Value = 6305.5
$ws.Range("I70").Value = 5619.788
$ws.Range("K70").Value = 5619.788
$ws.Range("M70").Value = -5349.788
$ws.Range("H73").Value = 6305.5
$ws.Range("I73").Value = 5619.788
$ws.Range("K73").Value = 5619.788
$ws.Range("M73").Value = -4683.788
$ws.Range("H80").Value = 3466.8333
$ws.Range("I80").Value = 3675.625
$ws.Range("J80").Value = 3049.25
$ws.Range("K80").Value = 3675.625
$ws.Range("L80").Value = 3049.25
$ws.Range("M80").Value = -2677.625
$ws.Range("N80").Value = -5045.25
$ws.Range("H83").Value = 3466.8333
$ws.Range("I83").Value = 3675.625
$ws.Range("J83").Value = 3049.25
$ws.Range("K83").Value = 18378.125
$ws.Range("L83").Value = 15246.25
$ws.Range("M83").Value = -13386.125
$ws.Range("N83").Value = -25230.25
$ws.Range("H132").Value = 2192.4
$ws.Range("I132").Value = 1160.8667
$ws.Range("J132").Value = 3430.24
$ws.Range("K132").Value = 3482.6001
$ws.Range("L132").Value = 10290.72
$ws.Range("M132").Value = -952.6001000000001
$ws.Range("N132").Value = -15350.72
$ws.Range("H137").Value = 30312
$ws.Range("J137").Value = 45780
$ws.Range("L137").Value = 45780
$ws.Range("N137").Value = -55980

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H16").Value = 1599.7142
$ws.Range("I16").Value = 1175
$ws.Range("J16").Value = 2166
$ws.Range("K16").Value = 1175
$ws.Range("L16").Value = 2166
$ws.Range("M16").Value = -1005
$ws.Range("N16").Value = -2506

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H108").Value = 30500
$ws.Range("J108").Value = 30500
$ws.Range("L108").Value = 30500
$ws.Range("N108").Value = -38180
$ws.Range("H113").Value = 457.3125
$ws.Range("I113").Value = 446.72726
$ws.Range("J113").Value = 480.6
$ws.Range("K113").Value = 1340.18178
$ws.Range("L113").Value = 1441.8
$ws.Range("M113").Value = 829.8182200000001
$ws.Range("N113").Value = -5781.8
